# Apply "Final version of technical document, presentation" pairing-list edit:
#  - copy the auxiliary name list (G3:G13) into J3:J13 and alphabetically sort it
#  - add the English transliteration of each (already-sorted) name in M3:M13
#  - add a 0/1 flag column in O3:O13
#  - color-scale the flag column
#  - leave the new "J3:M13" block selected, matching the author's final view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copy the existing G3:G13 list into J3:J13 (same Georgian names, still
#    in their original/unsorted order at this point).
$ws.Range("J3:J13").Value2 = $ws.Range("G3:G13").Value2

# 2) Sort J3:J13 ascending - this is what produces the <sortState>/<sortCondition>
#    block that Excel records after a manual sort of this range.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("J3")) | Out-Null
$sortRange = $ws.Range("J3:J13")
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlNo
$ws.Sort.MatchCase = $false
$ws.Sort.Orientation = [Microsoft.Office.Interop.Excel.XlSortOrientation]::xlSortColumns
$ws.Sort.Apply()

# 3) After the sort, J3:J13 is in a known, fixed (alphabetical) order - fill in
#    the matching English transliteration (M) and the 0/1 flag (O) per row.
$pairs = @(
    @{Row = 3;  M = "Kloyan Manvel";        O = 1},
    @{Row = 4;  M = "Korakhashvili Luka";   O = 1},
    @{Row = 5;  M = "Latsabidze Giorgi";    O = 0},
    @{Row = 6;  M = "Makandarashvili Gia";  O = 1},
    @{Row = 7;  M = "Mirzashvili Giorgi";   O = 0},
    @{Row = 8;  M = "Sabashvili Irakli";    O = 0},
    @{Row = 9;  M = "Uridia Daviti";        O = 1},
    @{Row = 10; M = "Kochladze Guram";      O = 0},
    @{Row = 11; M = "Shonia Saba";          O = 0},
    @{Row = 12; M = "Tsikelashvili Giorgi"; O = 1},
    @{Row = 13; M = "Chankvetadze Mariam";  O = 0}
)

foreach ($pair in $pairs) {
    $r = $pair.Row
    $ws.Cells.Item($r, 13).Value2 = $pair.M   # column M
    $ws.Cells.Item($r, 15).Value2 = $pair.O   # column O
}

# 4) Color-scale (red/yellow/green, Excel's default 3-color scale) the flag column.
$flagRange = $ws.Range("O3:O13")
$flagRange.FormatConditions.Delete()
$flagRange.FormatConditions.AddColorScale(3) | Out-Null

# 5) Leave the English-name column selected, with the view scrolled over so
#    column C is the leftmost visible column - matching the author's final state.
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("M3:M13").Select() | Out-Null
